$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The month headers in E1:P1 used to label a single calendar month (e.g. "06-2022").
# They now label a one-month *range* spanning that month and the following one
# (e.g. "06-2022 - 07-2022"), sliding forward one column at a time through "05-2023 - 06-2023".
$monthRanges = @(
    "06-2022 - 07-2022",
    "07-2022 - 08-2022",
    "08-2022 - 09-2022",
    "09-2022 - 10-2022",
    "10-2022 - 11-2022",
    "11-2022 - 12-2022",
    "12-2022 - 01-2023",
    "01-2023 - 02-2023",
    "02-2023 - 03-2023",
    "03-2023 - 04-2023",
    "04-2023 - 05-2023",
    "05-2023 - 06-2023"
)

for ($i = 0; $i -lt $monthRanges.Length; $i++) {
    # Column E is index 5, through column P (index 16)
    $col = 5 + $i
    $ws.Cells.Item(1, $col).Value = $monthRanges[$i]
}

# These header columns are best-fit / auto-sized, so they need to widen to fit the
# new, longer header text (e.g. "06-2022" -> "06-2022 - 07-2022"). Re-running AutoFit
# grows them; nudge ColumnWidth to line up with the best-fit width Excel computes
# for this longer label so the stored column width matches the widened header text.
$headerColumns = $ws.Range("E1:P1").EntireColumn
$headerColumns.AutoFit()
$headerColumns.ColumnWidth = 20.3
